$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5685
$ws.Range("B2").Value = 3667
$ws.Range("C2").Value = 3711
$ws.Range("D2").Value = 7207
$ws.Range("E2").Value = 6410
$ws.Range("F2").Value = 3771
